$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Content for the new "Stop gambling" article (row 19 / article #18)
$content = @"
h3: Stop gambling.
p.note: I seriously read something that is a life changing ayat. All ayats are life changing, but just focus on the context of this particular ayat. <span class=“lavendar”>Ayat 187.</span>
p: I am arranging a tech meetup tomorrow morning (Sunday 19 Jan 2020). While the visitors will be arriving, I planned to keep it musical until all of them arrived. I even tested it by playing some gym hardcore songs in the background. Plugged in the speakers and tested it all aloud. 
p: Now I am sitting on a desk, writing this blog. Acting all religious like nothing happened. Sins forgotten, came back to the Holy track. Tomorrow when the event comes, I go back pleasing my audience using all ‘feel good means’. 
p: I will be a hypocrite if I do not implement what I read today. Why? Because Allah (swt) mentions in this ayat:-
quote: These are the limits [set by] Allah, so do not approach them. <br> - Surah Baqarah verse 187
p: Allah (swt) mentions do not approach anything which is near the boundary of evil, where you might topple and fall on the other side. If you want to play football, play it in the center of the field. Do not go near the boundary and be a stuntsman. Come back to the center. These stunts might lead you to destruction. 
p: Thus, playing music loud out loud, will bring me closer to the boundary and I might topple. This is not a discussion of halal or haram. This is a discussion of staying sane and avoiding areas which might lead us to the other side of the boundary. According to this ayat, things that take us out of the playing ground, have to be avoided at all costs. <b>Approaching.. simply approaching the evil leads to destruction.</b>
h3: Some things that gradually pulls us down.
p.b-left: Do not smoke a ciggeratte because it is the gateway to other drugs. Not arguing if it is halal or haram. It just leads us out of the field.
p.b-left: Do not listen to music, because it leads to mental dullness. It makes us emotional, we start losing focus and go astray. 
p.b-left: Do not look for cheerleaders on the other side of the boundary, you will go out of the ground. Might never come back. Might die in such a state. Might lose the chance to succeed.
p.b-left: Start using word ‘Sister’ when you come across any other girl. I have used it, it really helped me. Satan defeated. My Nafs defeated. Sister defeated. Right where it lifted its head, we smashed it back to the ground. 
p: Islam is fun and peaceful. Anything that ridicules our religion. Anything that takes us away from the Core of the field, is a detraction. Following distractions does not make us any better then a 7 years old kid. <b>Our learnings and experiences are useless if we are stuck in distractions. </b>
h3: Also..
p: Stop looking for fatwas 
p.b-left: <b>“Cigerrete haram or halal”.</b> 
p.b-left: <b>“Music haram or halal”.</b> 
p: Ask yourself, if it affects your focus to your goals, it is just another distraction. Living our life saying, <b>“one day Hidayet will come and we will get back on track”</b>, is a big chance. May be we should get out of this Casino. No body gets rich in casinos. The house always wins. Satan always wins in the gamble of truth and lie. So why taking the chance. 
p: Let us get out of the Casino. Enjoy Las Vegas in open. Sit under clouds. Do some coding. Write some blogs. Smile and feel contended. Lets enjoy our lives. Lets stop taking chances..
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>.
"@

$ws.Range("A19").Value = 18
$ws.Range("B19").NumberFormat = "d-mmm-yy"
$ws.Range("B19").Value = Get-Date -Year 2020 -Month 1 -Day 18 -Hour 0 -Minute 0 -Second 0
$ws.Range("C19").Value = "Surah Baqarah, 183 - 188"
$ws.Range("D19").Value = $content
$ws.Range("E19").Value = "Qasim Ali"
$ws.Range("F19").Value = "Casinos, Boundary of Islam, Risk in lies, Satan attacks"

$ws.Rows.Item(19).RowHeight = 409.6

$ws.Application.ActiveWindow.ScrollRow = 19 | Out-Null
$ws.Range("F19").Select() | Out-Null
